$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.975.07"
$ws.Range("E2").Value = "  -4.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.494.01"
$ws.Range("E3").Value = "  -3.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "281.87"
$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3797"
$ws.Range("E7").Value = "  -4.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  -3.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.63"
$ws.Range("E9").Value = "  -2.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06867"
$ws.Range("E10").Value = "  -5.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.027"
$ws.Range("E11").Value = "  -4.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.515"
$ws.Range("E13").Value = "  -3.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.76"
$ws.Range("E14").Value = "  -5.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.502.18"
$ws.Range("E15").Value = "  -3.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.322"
$ws.Range("E16").Value = "  -4.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001058"
$ws.Range("E17").Value = "  -6.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06538"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.52"
$ws.Range("E19").Value = "  -2.47%  "

$ws.Range("E20").Value = "  +0.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.915"
$ws.Range("E21").Value = "  -5.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.95"
$ws.Range("E22").Value = "  -3.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.78"
$ws.Range("E23").Value = "  -4.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.344"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.980.20"
$ws.Range("E25").Value = "  -4.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.301"
$ws.Range("E26").Value = "  -4.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.65"
$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.80"
$ws.Range("E28").Value = "  -4.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.782"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.674.01"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "113.29"
$ws.Range("E31").Value = "  -4.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.786"
$ws.Range("E32").Value = "  -0.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9297"
$ws.Range("E33").Value = "  -3.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07888"
$ws.Range("E34").Value = "  -5.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.462"
$ws.Range("E35").Value = "  -7.40%  "

$ws.Range("E36").Value = "  -7.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.991"
$ws.Range("E37").Value = "  -2.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.97"
$ws.Range("E38").Value = "  +2.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05740"
$ws.Range("E39").Value = "  -4.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02101"
$ws.Range("E41").Value = "  -7.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.148"
$ws.Range("E42").Value = "  -5.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1955"
$ws.Range("E43").Value = "  -4.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5561"
$ws.Range("E44").Value = "  -4.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.71"
$ws.Range("E45").Value = "  -2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.641"
$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5380"
$ws.Range("E47").Value = "  -3.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.121"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.822"
$ws.Range("E49").Value = "  -4.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.30"
$ws.Range("E50").Value = "  -5.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06497"
$ws.Range("E51").Value = "  -4.68%  "
